$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp (A1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 14 de Mayo de 2020 a las 23:05"

# Country-name reshuffle for rows 190-198 (A column)
$ws.Cells.Item(190,1).Value = "Mauritania"
$ws.Cells.Item(191,1).Value = "Butan"
$ws.Cells.Item(192,1).Value = "Laos"
$ws.Cells.Item(193,1).Value = "Fiyi"
$ws.Cells.Item(196,1).Value = "Nueva Caledonia"
$ws.Cells.Item(197,1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(198,1).Value = "San Vicente y las Granadinas"

# Updated statistic cells (B..H)
$ws.Cells.Item(4,2).Value = 1451500
$ws.Cells.Item(4,3).Value = 21152
$ws.Cells.Item(4,4).Value = 316181
$ws.Cells.Item(4,5).Value = 1048720
$ws.Cells.Item(4,7).Value = 1402
$ws.Cells.Item(4,8).Value = 86599
$ws.Cells.Item(11,2).Value = 174909
$ws.Cells.Item(11,3).Value = 811
$ws.Cells.Item(11,5).Value = 16685
$ws.Cells.Item(11,7).Value = 63
$ws.Cells.Item(11,8).Value = 7924
$ws.Cells.Item(15,2).Value = 81997
$ws.Cells.Item(15,3).Value = 3942
$ws.Cells.Item(15,5).Value = 51379
$ws.Cells.Item(16,2).Value = 80604
$ws.Cells.Item(16,3).Value = 4298
$ws.Cells.Item(16,4).Value = 25151
$ws.Cells.Item(16,5).Value = 53186
$ws.Cells.Item(16,6).Value = 842
$ws.Cells.Item(16,7).Value = 98
$ws.Cells.Item(16,8).Value = 2267
$ws.Cells.Item(24,2).Value = 30502
$ws.Cells.Item(24,3).Value = 16
$ws.Cells.Item(24,5).Value = 24731
$ws.Cells.Item(24,6).Value = 196
$ws.Cells.Item(24,7).Value = 4
$ws.Cells.Item(24,8).Value = 2338
$ws.Cells.Item(30,4).Value = 5973
$ws.Cells.Item(30,5).Value = 20104
$ws.Cells.Item(63,2).Value = 5530
$ws.Cells.Item(63,3).Value = 122
$ws.Cells.Item(63,4).Value = 674
$ws.Cells.Item(63,5).Value = 4832
$ws.Cells.Item(75,2).Value = 2645
$ws.Cells.Item(75,3).Value = 33
$ws.Cells.Item(75,5).Value = 498
$ws.Cells.Item(101,6).Value = 13
$ws.Cells.Item(135,4).Value = 285
$ws.Cells.Item(135,5).Value = 24
$ws.Cells.Item(182,4).Value = 13
$ws.Cells.Item(182,5).Value = 20
$ws.Cells.Item(190,2).Value = 20
$ws.Cells.Item(190,3).Value = 5
$ws.Cells.Item(190,4).Value = 6
$ws.Cells.Item(190,5).Value = 12
$ws.Cells.Item(190,8).Value = 2
$ws.Cells.Item(191,3).Value = 8
$ws.Cells.Item(191,4).Value = 5
$ws.Cells.Item(191,5).Value = 14
$ws.Cells.Item(192,2).Value = 19
$ws.Cells.Item(192,5).Value = 5
$ws.Cells.Item(193,4).Value = 14
$ws.Cells.Item(193,5).Value = 4
$ws.Cells.Item(196,2).Value = 18
$ws.Cells.Item(196,4).Value = 18
$ws.Cells.Item(196,5).Value = 0
$ws.Cells.Item(197,4).Value = 0
$ws.Cells.Item(197,5).Value = 17
$ws.Cells.Item(198,2).Value = 17
$ws.Cells.Item(198,3).Value = 0
$ws.Cells.Item(198,4).Value = 12
$ws.Cells.Item(198,5).Value = 5
$ws.Cells.Item(198,8).Value = 0
